# Apply the "model and template with unit, description and enum" edit:
#  - row2 (type-hint row), column U (MolecularWeight) gains a unit annotation
#  - a new row 3 with French field descriptions for the first 7 columns is added

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update U2: #float -> #float,  unit:kb -----------------------------
$ws.Range("U2").Value = "#float,  unit:kb"

# --- Add new row 3: description row ------------------------------------
$row3 = @(
    "#Manipulateur",
    "#Desc:IdentifiantEchantillon",
    "#Date",
    "#ModeOderatoireLaboratoire",
    "#AppareilLogicielCritique",
    "#ProduitCritique",
    "#LieuStockageDonneesBrutes",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    "",
    ""
)

for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
